$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2027
$ws1.Range("F5").Value = 328
$ws1.Range("F7").Value = 94
$ws1.Range("F8").Value = 2056
$ws1.Range("F9").Value = 10466
$ws1.Range("F12").Value = 273
$ws1.Range("F14").Value = 404
$ws1.Range("F15").Value = 7349
$ws1.Range("F17").Value = 697
$ws1.Range("F18").Value = 179
$ws1.Range("F20").Value = 3285

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2027
$ws4.Range("F5").Value = 328
$ws4.Range("F8").Value = 94
$ws4.Range("F9").Value = 2056
$ws4.Range("F12").Value = 10466
$ws4.Range("F15").Value = 273
$ws4.Range("F17").Value = 404
$ws4.Range("F18").Value = 7349
$ws4.Range("F20").Value = 697
$ws4.Range("F21").Value = 179
$ws4.Range("F23").Value = 3285
